# Apply the edits described by the commit "run simulation and plot rna results"
$wb = $excel.ActiveWorkbook

# --- Model sheet: update the "Created"/"Updated" timestamp ---
# Both cells originally shared the same string value ("2019-08-16 18:41:28");
# the canonical edit updates that shared string in place, so both cells
# must be set to the new value to reproduce the effect.
$wsModel = $wb.Worksheets.Item("Model")
$wsModel.Range("B11").Value = "2019-08-19 12:43:06"
$wsModel.Range("B12").Value = "2019-08-19 12:43:06"

# --- Parameters sheet: update kinetic constant values ---
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("D11").Value = 0.003080654135821979
$wsParams.Range("D12").Value = 0.003080654135821979
$wsParams.Range("D13").Value = 0.003080654135821979
$wsParams.Range("D14").Value = 0.003080654135821979
$wsParams.Range("D18").Value = 0.001
$wsParams.Range("D19").Value = 0.001
$wsParams.Range("D23").Value = 0.001
$wsParams.Range("D24").Value = 0.001
$wsParams.Range("D25").Value = 0.001
$wsParams.Range("D26").Value = 0.001
$wsParams.Range("D28").Value = 0.001
$wsParams.Range("D29").Value = 0.001

# --- Initial species concentrations sheet: update mean concentrations ---
$wsConc = $wb.Worksheets.Item("Initial species concentrations")
$wsConc.Range("E2").Value = 30110.704285
$wsConc.Range("E3").Value = 30110.704285
$wsConc.Range("E5").Value = 30110.704285
$wsConc.Range("E6").Value = 30110.704285
$wsConc.Range("E8").Value = 30110.704285
$wsConc.Range("E9").Value = 30110.704285
$wsConc.Range("E11").Value = 1656088735.675
$wsConc.Range("E19").Value = 30110.704285
$wsConc.Range("E20").Value = 30110.704285
